$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Register")

# Update the Username value for the registered customer
$ws.Range("I2").Value = "Ars2001"

# Select the cell that was last edited, matching the saved selection state
$ws.Range("I2").Select()
